$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.143.07"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").Value = "1.559.49"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.09%  "

$ws.Range("E10").Value = "  -0.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "1.781.91"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("D13").Value = "1.558.65"
$ws.Range("E13").Value = "  -2.29%  "

$ws.Range("E14").Value = "  -2.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "27.148.70"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  -1.67%  "

$ws.Range("E20").Value = "  -2.02%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("E23").Value = "  -3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("E29").Value = "  -1.55%  "

$ws.Range("E30").Value = "  -2.09%  "

$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").Value = "1.380.62"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.946"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.72%  "

$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.811"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.511"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.96%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.94%  "

$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").Value = "1.694.41"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").Value = "0.0₇0982"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("E51").Value = "  +0.08%  "
